$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data rows (2-7) so the stale shared strings (FAPs, sCs, Bmp2, Bmpr1a, ECs)
# are dropped from the table and can be rebuilt in the desired order.
$ws.Range("A2:T10").Clear()

# Write the text columns column-by-column (A, then B, then C, then D) across all rows
# so the shared-string table is rebuilt in the order: ECs, FAPs, sCs, Bmp2, Bmpr1a.
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "ECs"
$ws.Range("A5").Value = "FAPs"
$ws.Range("A6").Value = "FAPs"
$ws.Range("A7").Value = "FAPs"
$ws.Range("A8").Value = "sCs"
$ws.Range("A9").Value = "sCs"
$ws.Range("A10").Value = "sCs"

$ws.Range("B2").Value = "Bmp2"
$ws.Range("B3").Value = "Bmp2"
$ws.Range("B4").Value = "Bmp2"
$ws.Range("B5").Value = "Bmp2"
$ws.Range("B6").Value = "Bmp2"
$ws.Range("B7").Value = "Bmp2"
$ws.Range("B8").Value = "Bmp2"
$ws.Range("B9").Value = "Bmp2"
$ws.Range("B10").Value = "Bmp2"

$ws.Range("C2").Value = "Bmpr1a"
$ws.Range("C3").Value = "Bmpr1a"
$ws.Range("C4").Value = "Bmpr1a"
$ws.Range("C5").Value = "Bmpr1a"
$ws.Range("C6").Value = "Bmpr1a"
$ws.Range("C7").Value = "Bmpr1a"
$ws.Range("C8").Value = "Bmpr1a"
$ws.Range("C9").Value = "Bmpr1a"
$ws.Range("C10").Value = "Bmpr1a"

$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "sCs"
$ws.Range("D5").Value = "ECs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "sCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("D10").Value = "sCs"

# Write the numeric columns for each row
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.445484
$ws.Range("H2").Value = 4.336452
$ws.Range("I2").Value = 0.1286708197254238
$ws.Range("J2").Value = 0.1286708197254238
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.918858999999999
$ws.Range("N2").Value = 29.756577
$ws.Range("O2").Value = 0.160764128269069
$ws.Range("P2").Value = 0.160764128269069
$ws.Range("Q2").Value = 14.337551982756
$ws.Range("R2").Value = 129.037967844804
$ws.Range("S2").Value = 0.02068565216682429
$ws.Range("T2").Value = 0.02068565216682429

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.445484
$ws.Range("H3").Value = 4.336452
$ws.Range("I3").Value = 0.1286708197254238
$ws.Range("J3").Value = 0.1286708197254238
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 29.20351433333333
$ws.Range("N3").Value = 87.61054300000001
$ws.Range("O3").Value = 0.47332838627826
$ws.Range("P3").Value = 0.4733283862782601
$ws.Range("Q3").Value = 42.213212712604
$ws.Range("R3").Value = 379.918914413436
$ws.Range("S3").Value = 0.06090355146173577
$ws.Range("T3").Value = 0.06090355146173577

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.445484
$ws.Range("H4").Value = 4.336452
$ws.Range("I4").Value = 0.1286708197254238
$ws.Range("J4").Value = 0.1286708197254238
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.575837
$ws.Range("N4").Value = 67.72751099999999
$ws.Range("O4").Value = 0.3659074854526709
$ws.Range("P4").Value = 0.3659074854526709
$ws.Range("Q4").Value = 32.63301117010799
$ws.Range("R4").Value = 293.6971005309719
$ws.Range("S4").Value = 0.04708161609686375
$ws.Range("T4").Value = 0.04708161609686376

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.292313
$ws.Range("H5").Value = 18.876939
$ws.Range("I5").Value = 0.5601148623429528
$ws.Range("J5").Value = 0.5601148623429528
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.918858999999999
$ws.Range("N5").Value = 29.756577
$ws.Range("O5").Value = 0.160764128269069
$ws.Range("P5").Value = 0.160764128269069
$ws.Range("Q5").Value = 62.412565430867
$ws.Range("R5").Value = 561.713088877803
$ws.Range("S5").Value = 0.09004637757511441
$ws.Range("T5").Value = 0.09004637757511441

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.292313
$ws.Range("H6").Value = 18.876939
$ws.Range("I6").Value = 0.5601148623429528
$ws.Range("J6").Value = 0.5601148623429528
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 29.20351433333333
$ws.Range("N6").Value = 87.61054300000001
$ws.Range("O6").Value = 0.47332838627826
$ws.Range("P6").Value = 0.4733283862782601
$ws.Range("Q6").Value = 183.7576528853197
$ws.Range("R6").Value = 1653.818875967877
$ws.Range("S6").Value = 0.2651182639232596
$ws.Range("T6").Value = 0.2651182639232597

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.292313
$ws.Range("H7").Value = 18.876939
$ws.Range("I7").Value = 0.5601148623429528
$ws.Range("J7").Value = 0.5601148623429528
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.575837
$ws.Range("N7").Value = 67.72751099999999
$ws.Range("O7").Value = 0.3659074854526709
$ws.Range("P7").Value = 0.3659074854526709
$ws.Range("Q7").Value = 142.054232640981
$ws.Range("R7").Value = 1278.488093768829
$ws.Range("S7").Value = 0.2049502208445788
$ws.Range("T7").Value = 0.2049502208445788

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.496172
$ws.Range("H8").Value = 10.488516
$ws.Range("I8").Value = 0.3112143179316233
$ws.Range("J8").Value = 0.3112143179316232
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 9.918858999999999
$ws.Range("N8").Value = 29.756577
$ws.Range("O8").Value = 0.160764128269069
$ws.Range("P8").Value = 0.160764128269069
$ws.Range("Q8").Value = 34.678037107748
$ws.Range("R8").Value = 312.102333969732
$ws.Range("S8").Value = 0.05003209852713032
$ws.Range("T8").Value = 0.05003209852713031

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.496172
$ws.Range("H9").Value = 10.488516
$ws.Range("I9").Value = 0.3112143179316233
$ws.Range("J9").Value = 0.3112143179316232
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 29.20351433333333
$ws.Range("N9").Value = 87.61054300000001
$ws.Range("O9").Value = 0.47332838627826
$ws.Range("P9").Value = 0.4733283862782601
$ws.Range("Q9").Value = 102.1005091137987
$ws.Range("R9").Value = 918.9045820241881
$ws.Range("S9").Value = 0.1473065708932646
$ws.Range("T9").Value = 0.1473065708932646

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.496172
$ws.Range("H10").Value = 10.488516
$ws.Range("I10").Value = 0.3112143179316233
$ws.Range("J10").Value = 0.3112143179316232
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 22.575837
$ws.Range("N10").Value = 67.72751099999999
$ws.Range("O10").Value = 0.3659074854526709
$ws.Range("P10").Value = 0.3659074854526709
$ws.Range("Q10").Value = 78.92900919596399
$ws.Range("R10").Value = 710.361082763676
$ws.Range("S10").Value = 0.1138756485112283
$ws.Range("T10").Value = 0.1138756485112283

